$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) row data: title, timestamp, historical distance, time bucket, uri
# Rows shuffled per the commit "added one json for time bucket analysis"

$ws.Range("A2").Value = "Passengers diverted to ONT try to wind way back to their destinations"
$ws.Range("B2").Value = "2013-11-01T00:00:00UTC"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "day_0"
$ws.Range("E2").Value = "https://web.archive.org/web/20131102211206/http://www.pe.com/local-news/local-news-headlines/20131101-some-flights-diverted-to-ont-after-lax-shooting.ece"

$ws.Range("A3").Value = "LA airport attack: Paul Ciancia pleads not guilty"
$ws.Range("B3").Value = "2013-12-26T19:44:56UTC"
$ws.Range("C3").Value = 55
$ws.Range("D3").Value = "day_31_beyond"
$ws.Range("E3").Value = "https://www.bbc.co.uk/news/world-us-canada-25523003"

$ws.Range("A4").Value = "Suspected LAX gunman had his targets clearly in mind"
$ws.Range("B4").Value = "2013-11-03T03:00:00UTC"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "day_2_to_30"
$ws.Range("E4").Value = "https://web.archive.org/web/20131104213951/http://www.latimes.com/local/la-me-1103-lax-shooting-20131103,0,7342159.story#axzz2jmqKVQgo"

$ws.Range("A5").Value = "H.R. 4802, Gerardo Hernandez Airport Security Act of 2014"
$ws.Range("B5").Value = "2014-06-24T00:00:00UTC"
$ws.Range("C5").Value = 235
$ws.Range("D5").Value = "day_31_beyond"
$ws.Range("E5").Value = "https://web.archive.org/web/20140728020009/http://www.cbo.gov/publication/45476"

$ws.Range("A6").Value = "LAX shooting: Notice of intent to seek death penalty"
$ws.Range("B6").Value = "2015-01-02T00:00:00UTC"
$ws.Range("C6").Value = 427
$ws.Range("D6").Value = "day_31_beyond"
$ws.Range("E6").Value = "http://documents.latimes.com/lax-shooting-notice-intent-seek-death-penalty/"

$ws.Range("A7").Value = "AP: LAX shooter opened fire while officers on break"
$ws.Range("B7").Value = "2014-01-22T15:39:00UTC"
$ws.Range("C7").Value = 82
$ws.Range("D7").Value = "day_31_beyond"
$ws.Range("E7").Value = "http://www.cbsnews.com/news/ap-lax-shooter-opened-fire-while-officers-on-break/"

$ws.Range("A8").Value = "Man who killed TSA officer at LAX in 2013 pleads guilty and avoids death penalty"
$ws.Range("B8").Value = "2016-09-06T00:00:00UTC"
$ws.Range("C8").Value = 1040
$ws.Range("D8").Value = "day_31_beyond"
$ws.Range("E8").Value = "https://web.archive.org/web/20160913111152/http://www.latimes.com/local/lanow/la-me-lax-shooter-20160906-snap-story.html"

$ws.Range("A9").Value = "TSA agent shot, killed at LAX, suspect in critical condition"
$ws.Range("B9").Value = "2013-11-01T19:21:35UTC"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "day_0"
$ws.Range("E9").Value = "http://www.upi.com/Top_News/US/2013/11/01/Reports-Suspect-dies-TSA-agent-shot-killed-at-LAX-shooting/UPI-80811383326473/?spt=rln&or=2"

$ws.Range("A10").Value = "Horrific LAX Shooting Scene"
$ws.Range("B10").Value = "2013-11-02T00:59:09UTC"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "day_1"
$ws.Range("E10").Value = "https://web.archive.org/web/20131222160450/http://www.cnn.com/video/data/2.0/video/bestoftv/2013/11/02/bb-lax-shooting-scene-tim-daly.cnn.html"

$ws.Range("A11").Value = "LAX Shooter Paul Anthony Ciancia's Personal Details Emerge"
$ws.Range("B11").Value = "2013-11-01T17:43:26UTC"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "day_0"
$ws.Range("E11").Value = "http://www.mediaite.com/online/lax-shooter-paul-anthony-ciancias-personal-details-emerge/"

$ws.Range("A12").Value = "LAX Shooting Rampage Trial Slated for 2016"
$ws.Range("B12").Value = "2017-06-30T15:47:22UTC"
$ws.Range("C12").Value = 1337
$ws.Range("D12").Value = "day_31_beyond"
$ws.Range("E12").Value = "http://www.nbclosangeles.com/news/local/LAX-TSA-Shooting-Rampage-Paul-Ciancia-Death-Penalty-287505621.html"

$ws.Range("A13").Value = "TSA agent shot at Los Angeles airport died in two to five minutes"
$ws.Range("B13").Value = "2013-11-20T21:20:00UTC"
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = "day_2_to_30"
$ws.Range("E13").Value = "https://web.archive.org/web/20131121084229/http://www.latimes.com/local/la-me-1121-tsa-agent-20131121,0,3031589.story"

# The hyperlink targets keep the same rId-per-row association as before the
# edit (the relationship part itself is untouched) - only the URL "location"
# fragment that used to sit on the E5 link now sits on the E4 link instead.
# Mutate the existing Hyperlink objects in place (rather than delete/re-add)
# so the <hyperlink ref="..."/> entries keep their original single-cell form.
$idx = 0
foreach ($hl in $ws.Hyperlinks) {
    $idx = $idx + 1
    if ($idx -eq 3) {
        $hl.SubAddress = "axzz2jmqKVQgo"
    }
    if ($idx -eq 4) {
        $hl.SubAddress = ""
    }
}
